$d = $word.ActiveDocument

# --- Change 1: letter date "September 19, 2025" -> "September 21, 2025" ---
$d.Content.Find.Execute("September 19, 2025", $true, $false, $false, $false, $false, $true, 1, $false, "September 21, 2025", 2)

# --- Change 2: split the mailing-address paragraph into three paragraphs ---
#   "nan PO Box 2911, Santa Clara CA 95055"
# becomes
#   "nan PO Box 2911"
#   "Santa Clara, CA 95055"
#   "" (new empty paragraph)
$full = $d.Content
$found = $full.Find.Execute("nan PO Box 2911, Santa Clara CA 95055")

$prefix = "nan PO Box 2911"
$splitPoint = $full.Start + $prefix.Length

# Remove the ", " that used to join the two halves of the address.
$commaRange = $d.Range($splitPoint, $splitPoint + 2)
$commaRange.Text = ""

# Break the paragraph right after "nan PO Box 2911".
$breakRange = $d.Range($splitPoint, $splitPoint)
$breakRange.InsertParagraphAfter()

# Re-find the now-standalone city/state/zip text and add the missing comma
# after "Santa Clara".
$full2 = $d.Content
$found2 = $full2.Find.Execute("Santa Clara CA 95055")
$santaClaraLen = "Santa Clara".Length
$commaInsertPoint = $full2.Start + $santaClaraLen
$insertRange = $d.Range($commaInsertPoint, $commaInsertPoint)
$insertRange.InsertAfter(",")

# Add the trailing empty paragraph after "Santa Clara, CA 95055".
$full3 = $d.Content
$found3 = $full3.Find.Execute("Santa Clara, CA 95055")
$endPoint = $full3.End
$finalBreakRange = $d.Range($endPoint, $endPoint)
$finalBreakRange.InsertParagraphAfter()

# --- Change 3: remove the two empty paragraphs that used to sit right
# after "Board of Directors" (a No Spacing-styled empty paragraph followed
# by a Title-styled empty paragraph); a third, also-empty Title paragraph
# stays in place. ---
$boardRange = $d.Content
$boardFound = $boardRange.Find.Execute("Board of Directors")
$boardStart = $boardRange.Start

$boardPara = $d.Range($boardStart, $boardStart).Paragraphs(1)
$firstEmpty = $boardPara.Next()
$d.Range($firstEmpty.Range.Start, $firstEmpty.Range.End).Delete()

$boardPara2 = $d.Range($boardStart, $boardStart).Paragraphs(1)
$secondEmpty = $boardPara2.Next()
$d.Range($secondEmpty.Range.Start, $secondEmpty.Range.End).Delete()
